$d = $word.ActiveDocument

function Find-Range($searchText) {
    # Search the whole document content for $searchText and return a Range
    # positioned exactly over the match (does not modify the document).
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, `
                              $true, 1, $false, $null, 0)
    if (-not $found) {
        throw "Could not find: $searchText"
    }
    return $r
}

# 1) Merge the "{x3}" run and the "{g4}ommentscay{/g5}" run (both inside the
#    comment range, directly adjacent) into a single run reading
#    "{x3}ommentscay". Replacing the text of the second run causes it to
#    coalesce with the identically-formatted preceding "{x3}" run, and both
#    runs stay safely inside the commentRangeStart/commentRangeEnd bounds.
$rComments = Find-Range("{g4}ommentscay{/g5}")
$rComments.Text = "ommentscay"

# 2) Turn the "{x6}" marker run (just after commentRangeEnd) into "{x4}".
#    This run's formatting matches the run before commentRangeEnd, so
#    replacing its whole text would coalesce it backwards across the
#    commentRangeEnd marker. Editing only the single digit character inside
#    the run avoids that whole-run-merge and keeps commentRangeEnd in place.
$rX6 = Find-Range("{x6}")
$digitPos = $rX6.Start + 2
$d.Range($digitPos, $digitPos + 1).Text = "4"

# 3) Turn the "{x7}" marker run into "{x5}", again via a single-character
#    edit so it keeps standing on its own (it will be merged with the
#    following run in the next step instead).
$rX7 = Find-Range("{x7}")
$digitPos2 = $rX7.Start + 2
$d.Range($digitPos2, $digitPos2 + 1).Text = "5"

# 4) Merge the "{g8}.{/g9}" run into the preceding "{x5}" run, producing a
#    single run reading "{x5}.".
$rDot = Find-Range("{g8}.{/g9}")
$rDot.Text = "."
